# Update the lattice-multiplication exercise table: each cell shows a
# problem "AB x CD" followed by 4 lines built from the digits of the
# two factors. This script rewrites each cell's text to the new
# exercise while preserving the existing run formatting (sz=32) and
# the <w:br/> line-break structure, by writing a string that uses the
# manual-line-break character (Chr(11) / vertical tab) in place of the
# original <w:br/> elements.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11

# Each entry: row, col, and the new lines for that cell.
# A $null line3 means the cell (per the source edit) ends up missing
# that digit row, leaving an extra blank line instead.
$cells = @(
    @{ Row=1; Col=1; L1="76 x 42"; L2="  4    2"; L3="7|    |"; L4="6|    |" },
    @{ Row=1; Col=2; L1="68 x 54"; L2="  5    4"; L3="6|    |"; L4="8|    |" },
    @{ Row=1; Col=3; L1="81 x 62"; L2="  6    2"; L3="8|    |"; L4="1|    |" },

    @{ Row=2; Col=1; L1="86 x 14"; L2="  1    4"; L3="8|    |"; L4="6|    |" },
    @{ Row=2; Col=2; L1="49 x 15"; L2="  1    5"; L3="4|    |"; L4="9|    |" },
    @{ Row=2; Col=3; L1="22 x 72"; L2="  7    2"; L3=$null;    L4="2|    |" },

    @{ Row=3; Col=1; L1="43 x 11"; L2="  1    1"; L3="4|    |"; L4="3|    |" },
    @{ Row=3; Col=2; L1="70 x 82"; L2="  8    2"; L3="7|    |"; L4="0|    |" },
    @{ Row=3; Col=3; L1="89 x 31"; L2="  3    1"; L3="8|    |"; L4="9|    |" },

    @{ Row=4; Col=1; L1="33 x 26"; L2="  2    6"; L3="3|    |"; L4="3|    |" },
    @{ Row=4; Col=2; L1="27 x 33"; L2="  3    3"; L3="2|    |"; L4="7|    |" },
    @{ Row=4; Col=3; L1="80 x 84"; L2="  8    4"; L3="8|    |"; L4="0|    |" },

    @{ Row=5; Col=1; L1="49 x 43"; L2="  4    3"; L3="4|    |"; L4="9|    |" },
    @{ Row=5; Col=2; L1="14 x 61"; L2="  6    1"; L3="1|    |"; L4="4|    |" },
    @{ Row=5; Col=3; L1="65 x 84"; L2="  8    4"; L3="6|    |"; L4="5|    |" }
)

foreach ($def in $cells) {
    $cell = $t.Cell($def.Row, $def.Col)
    $rng = $cell.Range

    $parts = @($def.L1, $def.L2, "  ----")
    if ($def.L3 -ne $null) {
        $parts += $def.L3
    } else {
        $parts += ""
    }
    $parts += $def.L4

    $newText = [string]::Join($vt, $parts)
    $rng.Text = $newText
}
